$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (row 11) Right count: 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row (row 12) Right total: 27 -> 45
$ws.Range("B12").Value = 45

# Update the correct/total marks text: 15/84 -> 45/140
$ws.Range("E12").Value = "45/140"
